$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.954.11"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.232.14"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.23%  "
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "2.572.10"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "2.254.11"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.810"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.52%  "
$ws.Range("D18").Value = "43.843.64"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "0.0₃0949"
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.69%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.43%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("E38").Value = "  -9.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0295"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.54%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "1.735.24"
$ws.Range("E44").Value = "  -4.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "84.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.184"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.67%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.67%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.42%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.92%  "
